# Update the two-digit-division answer table.
# Cells are addressed by (row, column) rather than by text search because
# some of the old values (e.g. "42÷3=14, 0") occur more than once in the
# table but map to different new values, so a blind Find/Replace-All would
# be unsafe.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "27÷8=3, 3"
$t.Cell(1, 2).Range.Text = "49÷9=5, 4"
$t.Cell(1, 3).Range.Text = "78÷6=13, 0"
$t.Cell(1, 4).Range.Text = "52÷3=17, 1"
$t.Cell(1, 5).Range.Text = "97÷7=13, 6"

# Row 5
$t.Cell(5, 1).Range.Text = "21÷7=3, 0"
$t.Cell(5, 2).Range.Text = "99÷6=16, 3"
$t.Cell(5, 3).Range.Text = "38÷4=9, 2"
$t.Cell(5, 4).Range.Text = "98÷7=14, 0"
$t.Cell(5, 5).Range.Text = "24÷9=2, 6"

# Row 9 (Cell(9,1) "31÷7=4, 3" is unchanged)
$t.Cell(9, 2).Range.Text = "50÷8=6, 2"
$t.Cell(9, 3).Range.Text = "97÷5=19, 2"
$t.Cell(9, 4).Range.Text = "66÷2=33, 0"
$t.Cell(9, 5).Range.Text = "51÷5=10, 1"

# Row 13
$t.Cell(13, 1).Range.Text = "45÷2=22, 1"
$t.Cell(13, 2).Range.Text = "80÷5=16, 0"
$t.Cell(13, 3).Range.Text = "96÷4=24, 0"
$t.Cell(13, 4).Range.Text = "32÷4=8, 0"
$t.Cell(13, 5).Range.Text = "21÷3=7, 0"

# Row 17
$t.Cell(17, 1).Range.Text = "34÷4=8, 2"
$t.Cell(17, 2).Range.Text = "24÷3=8, 0"
$t.Cell(17, 3).Range.Text = "29÷8=3, 5"
$t.Cell(17, 4).Range.Text = "84÷9=9, 3"
$t.Cell(17, 5).Range.Text = "99÷7=14, 1"
